# The workbook originally has two sheets, in this tab order:
#   1) hotel_info   - one header row + one data row (hotel-level info)
#   2) review_info  - header row only (review-level info), no data rows
#
# Target state (per the authoritative diff):
#   1) review_info  - unchanged headers, still no data rows, now FIRST
#   2) hotel_info   - SECOND, with a new "State" column inserted right
#                     after "Hotel_Name" (so columns become:
#                     STR, Hotel_Name, State, City, Zip, TA_ReviewURL,
#                     Tripadvisor_Hotel_Name, English_Reviews_num,
#                     Local_Rank, Total_Reviews_num), and the data row
#                     gets "Louisiana" in the new State column.

$wb = $excel.ActiveWorkbook

$hotelSheet  = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")

# 1. Insert a new column C ("State") into hotel_info, shifting
#    City/Zip/TA_ReviewURL/... one column to the right, and populate it.
$hotelSheet.Columns.Item(3).Insert()
$hotelSheet.Range("C1").Value = "State"
$hotelSheet.Range("C2").Value = "Louisiana"

# 2. Reorder the sheet tabs so review_info comes before hotel_info.
$reviewSheet.Move($hotelSheet)
